$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 39
$ws1.Range("F3").Value = 107
$ws1.Range("F4").Value = 1516
$ws1.Range("F5").Value = 217
$ws1.Range("F7").Value = 428
$ws1.Range("F8").Value = 9939
$ws1.Range("F10").Value = 123
$ws1.Range("F11").Value = 242
$ws1.Range("F12").Value = 188
$ws1.Range("F13").Value = 376
$ws1.Range("F14").Value = 6873
$ws1.Range("F16").Value = 638
$ws1.Range("F17").Value = 52

# Sheet "全部类型" (all types, combined list)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 39
$ws4.Range("F3").Value = 107
$ws4.Range("F4").Value = 1516
$ws4.Range("F5").Value = 217
$ws4.Range("F8").Value = 428
$ws4.Range("F11").Value = 9939
$ws4.Range("F13").Value = 123
$ws4.Range("F14").Value = 242
$ws4.Range("F15").Value = 188
$ws4.Range("F16").Value = 376
$ws4.Range("F17").Value = 6873
$ws4.Range("F19").Value = 638
$ws4.Range("F20").Value = 52
